# Add a new "2023" column (column U) to the statistics table, copying the
# formatting of the existing last column (T, "2022") and filling in the
# new year's values for each indicator row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column T (rows 4-14, the header/year row through the last data row)
# into column U so the new column inherits the same number formats, fonts,
# borders and alignment as the rest of the table.
$ws.Range("T4:T14").Copy($ws.Range("U4:U14"))

# Overwrite the copied values with the real 2023 data for each row.
$ws.Range("U4").Value2 = 2023
$ws.Range("U5").Value2 = 3.3
$ws.Range("U6").Value2 = 1
$ws.Range("U7").Value2 = 1.6
$ws.Range("U8").Value2 = 9.1999999999999993
$ws.Range("U9").Value2 = 6.1
$ws.Range("U10").Value2 = 1.5
$ws.Range("U11").Value2 = 4
$ws.Range("U12").Value2 = 4.4000000000000004
$ws.Range("U13").Value2 = 4.7
$ws.Range("U14").Value2 = 0.5

# Restore the selection to match the final state of the sheet (B1).
[void]$ws.Range("B1").Select()
